# Update the date line.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-09-03 Wednesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-09-04 Thursday", 2)

# Update the practice-problem table in place, cell by cell, so that
# duplicate / unchanged values (e.g. "31÷2=" staying "31÷2=", or two
# different cells both becoming "88÷5=") are handled unambiguously.
$t = $d.Tables.Item(1)

$values = @(
    @{Row=1;  Col=1; Text="94÷8="},
    @{Row=1;  Col=2; Text="68÷5="},
    @{Row=1;  Col=3; Text="88÷5="},
    @{Row=1;  Col=4; Text="72÷2="},
    @{Row=1;  Col=5; Text="22÷9="},

    @{Row=5;  Col=1; Text="66÷7="},
    @{Row=5;  Col=2; Text="62÷9="},
    @{Row=5;  Col=3; Text="77÷9="},
    @{Row=5;  Col=4; Text="19÷3="},
    @{Row=5;  Col=5; Text="88÷5="},

    @{Row=9;  Col=1; Text="64÷6="},
    @{Row=9;  Col=2; Text="41÷7="},
    @{Row=9;  Col=3; Text="17÷9="},
    @{Row=9;  Col=4; Text="40÷4="},
    @{Row=9;  Col=5; Text="57÷9="},

    @{Row=13; Col=1; Text="17÷7="},
    @{Row=13; Col=2; Text="33÷8="},
    @{Row=13; Col=3; Text="49÷4="},
    @{Row=13; Col=4; Text="35÷8="},
    @{Row=13; Col=5; Text="82÷2="},

    @{Row=17; Col=1; Text="31÷2="},
    @{Row=17; Col=2; Text="97÷4="},
    @{Row=17; Col=3; Text="65÷6="},
    @{Row=17; Col=4; Text="75÷2="},
    @{Row=17; Col=5; Text="37÷8="}
)

foreach ($v in $values) {
    $t.Cell($v.Row, $v.Col).Range.Text = $v.Text
}
